$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the gramStart/gramEnd proofErr wrapping around "( Will"
# and merge the "( Will be pinpointed by days/week available )" text into a
# single run (right after the preceding manual line break), instead of the
# three runs (with proofErr markers) that currently hold it.
# ---------------------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("( Will be pinpointed by days/week available )", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the '( Will ... )' text to edit"
}
$start1 = $find1.Start
$end1 = $find1.End

# Delete the manual line break immediately before the text together with the
# "( Will ... )" run(s) (this also drops the gramStart/gramEnd proofErr
# marks, since both are fully inside the deleted range).
$targetRange = $d.Range($start1 - 1, $end1)
$targetRange.Delete()

# Re-insert the line break + text as one chunk at the same spot. Because the
# insertion point sits right against the still-present run boundary, the
# engine folds the new break+text into that run instead of spawning fresh
# proofErr-wrapped runs.
$ins = $d.Range($start1 - 1, $start1 - 1)
$newChunk = [char]11 + "( Will be pinpointed by days/week available )"
$ins.InsertAfter($newChunk)

# ---------------------------------------------------------------------------
# Change 2: drop the "Feedback form" and "Admin dashboard" bullet lines
# (each introduced by a manual line break) from the Functional Requirements
# paragraph, keeping the list ending at "Diet recommendations ...".
# ---------------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("- Feedback form: Allows users to send feedback or questions.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the 'Feedback form' text to edit"
}
$feedbackStart = $find2.Start

$find3 = $d.Content
$found3 = $find3.Find.Execute("- Admin dashboard: Enables admin users to manage profiles and plans.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not locate the 'Admin dashboard' text to edit"
}
$adminEnd = $find3.End

# Include the manual line break right before "- Feedback form" so the whole
# two-bullet tail (breaks included) is removed in one shot.
$removeRange = $d.Range($feedbackStart - 1, $adminEnd)
$removeRange.Delete()
